# Applies the cryptos.xlsx price/volume update described in the commit
# "Updated cryptos list on Thu Aug 22 19:37:07 UTC 2024 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.125.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.18%  "
$ws.Range("D3").Value = "'2.597.12"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.56%  "
$ws.Range("D5").Value = "'583.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.74%  "
$ws.Range("D6").Value = "'142.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D8").Value = "'0.599"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.64%  "
$ws.Range("D9").Value = "'6.54"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.87%  "
$ws.Range("E10").Value = "  -1.03%  "
$ws.Range("E11").Value = "  -1.90%  "
$ws.Range("D12").Value = "'0.371"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.98%  "
$ws.Range("D13").Value = "'3.059.07"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.47%  "
$ws.Range("D14").Value = "'24.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.83%  "
$ws.Range("D15").Value = "'60.129.11"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.02%  "
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("D17").Value = "'2.601.91"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'11.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.21%  "
$ws.Range("E19").Value = "  -1.47%  "
$ws.Range("D20").Value = "'345.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.90%  "
$ws.Range("E21").Value = "  -1.74%  "
$ws.Range("E23").Value = "  +2.61%  "
$ws.Range("D24").Value = "'63.62"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.30%  "
$ws.Range("E25").Value = "  -0.54%  "
$ws.Range("D26").Value = "'0.159"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.61%  "
$ws.Range("D27").Value = "'7.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.06%  "
$ws.Range("E28").Value = "  +7.01%  "
$ws.Range("D29").Value = "'0.0₃0796"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("D30").Value = "'6.39"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.85%  "
$ws.Range("D31").Value = "'0.998"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Value = "'166.84"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.26%  "
$ws.Range("D33").Value = "'19.40"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.02%  "
$ws.Range("D34").Value = "'1.31"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +9.46%  "
$ws.Range("E35").Value = "  +0.73%  "
$ws.Range("D36").Value = "'0.980"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.75%  "
$ws.Range("E37").Value = "  +3.17%  "
$ws.Range("D38").Value = "'38.15"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.20%  "
$ws.Range("D39").Value = "'313.31"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.27%  "
$ws.Range("D40").Value = "'3.87"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.47%  "
$ws.Range("D41").Value = "'0.842"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.65%  "
$ws.Range("D42").Value = "'135.78"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.87%  "
$ws.Range("E43").Value = "  +0.86%  "
$ws.Range("D44").Value = "'1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.56%  "
$ws.Range("D45").Value = "'19.81"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.40%  "
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'4.97"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.89%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0242"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("D50").Value = "'19.85"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.03%  "
$ws.Range("D51").Value = "'10.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.49%  "
